$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B21").Value = "CN"
$ws.Range("C21").Value = "Aggressive Cow"
$ws.Range("D21").Value = "Java"
$ws.Range("E21").Value = "Medium"
$ws.Range("F21").Value = "To find min distance which is maximum"

# Match the "Medium" fill formatting used in other rows (e.g. E17) for the Level column
$ws.Range("E17").Copy()
$ws.Range("E21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("F21").Select()
